$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 651.3125
$ws.Range("I43").Value = 624.4286
$ws.Range("J43").Value = 672.2222
$ws.Range("K43").Value = 624.4286
$ws.Range("L43").Value = 672.2222
$ws.Range("M43").Value = -555.4286
$ws.Range("N43").Value = -810.2222

$ws.Range("I76").Value = 4118318.8
$ws.Range("J76").Value = 3966.6667
$ws.Range("K76").Value = 4118318.8
$ws.Range("L76").Value = 3966.6667
$ws.Range("M76").Value = -4118003.8

$ws.Range("I79").Value = 4118318.8
$ws.Range("J79").Value = 3966.6667
$ws.Range("K79").Value = 4118318.8
$ws.Range("L79").Value = 3966.6667
$ws.Range("M79").Value = -4117226.8

$ws.Range("H116").Value = 2440.625
$ws.Range("I116").Value = 2715.8333
$ws.Range("J116").Value = 1615
$ws.Range("K116").Value = 2715.8333
$ws.Range("L116").Value = 1615
$ws.Range("M116").Value = 726.1667000000002
$ws.Range("N116").Value = -8499

$ws.Range("H129").Value = 1148.1364
$ws.Range("I129").Value = 268.2
$ws.Range("J129").Value = 1406.9412
$ws.Range("K129").Value = 804.5999999999999
$ws.Range("L129").Value = 4220.8236
$ws.Range("M129").Value = 4195.4
$ws.Range("N129").Value = -14220.8236

$ws.Range("H137").Value = 38463084
$ws.Range("I137").Value = 47620216
$ws.Range("J137").Value = 3136.6
$ws.Range("K137").Value = 142860648
$ws.Range("L137").Value = 9409.799999999999
$ws.Range("M137").Value = -142858098
$ws.Range("N137").Value = -14509.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6155.5386
$ws.Range("I74").Value = 1154
$ws.Range("J74").Value = 11157.077
$ws.Range("K74").Value = 1154
$ws.Range("L74").Value = 11157.077
$ws.Range("M74").Value = -280
$ws.Range("N74").Value = -12905.077

$ws.Range("H77").Value = 6155.5386
$ws.Range("I77").Value = 1154
$ws.Range("J77").Value = 11157.077
$ws.Range("K77").Value = 5770
$ws.Range("L77").Value = 55785.38499999999
$ws.Range("M77").Value = -1402
$ws.Range("N77").Value = -64521.38499999999

$ws.Range("H122").Value = 2866.6667
$ws.Range("I122").Value = 2136.3635
$ws.Range("J122").Value = 4875
$ws.Range("K122").Value = 6409.0905
$ws.Range("L122").Value = 14625
$ws.Range("M122").Value = -3959.0905

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1106.36
$ws.Range("I20").Value = 841.8125
$ws.Range("J20").Value = 1576.6666
$ws.Range("K20").Value = 841.8125
$ws.Range("L20").Value = 1576.6666
$ws.Range("M20").Value = -594.8125
$ws.Range("N20").Value = -2070.6666

$ws.Range("N49").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 0

$ws.Range("H94").Value = 1268.48
$ws.Range("I94").Value = 1255.091
$ws.Range("J94").Value = 1366.6666
$ws.Range("K94").Value = 1255.091
$ws.Range("L94").Value = 1366.6666
$ws.Range("M94").Value = -804.0909999999999

$ws.Range("H134").Value = 4337.6
$ws.Range("I134").Value = 2243.5
$ws.Range("J134").Value = 7478.75
$ws.Range("K134").Value = 6730.5
$ws.Range("L134").Value = 22436.25
$ws.Range("M134").Value = -4195.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2015.6086
$ws.Range("I132").Value = 1146
$ws.Range("J132").Value = 4479.5
$ws.Range("K132").Value = 3438
$ws.Range("L132").Value = 13438.5
$ws.Range("M132").Value = -908

$ws.Range("H134").Value = 2882.4644
$ws.Range("I134").Value = 1441.3334
$ws.Range("J134").Value = 5476.5
$ws.Range("K134").Value = 4324.0002
$ws.Range("L134").Value = 16429.5
$ws.Range("M134").Value = -1789.0002
$ws.Range("N134").Value = -21499.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2256.375
$ws.Range("I34").Value = 2352
$ws.Range("J34").Value = 2242.7144
$ws.Range("K34").Value = 7056
$ws.Range("L34").Value = 6728.1432
$ws.Range("M34").Value = -6972
$ws.Range("N34").Value = -6896.1432

$ws.Range("H39").Value = 8559.200000000001
$ws.Range("I39").Value = 480
$ws.Range("J39").Value = 8895.833000000001
$ws.Range("K39").Value = 1440
$ws.Range("L39").Value = 26687.499
$ws.Range("M39").Value = -1146
$ws.Range("N39").Value = -27275.499

$ws.Range("M55").ClearContents()
$ws.Range("H55").Value = 3146.1538
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3146.1538
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 9438.4614
$ws.Range("N55").Value = -9792.4614

$ws.Range("H107").Value = 392.22726
$ws.Range("I107").Value = 431.91306
$ws.Range("J107").Value = 348.7619
$ws.Range("K107").Value = 1295.73918
$ws.Range("L107").Value = 1046.2857
$ws.Range("M107").Value = 624.26082
$ws.Range("N107").Value = -4886.2857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 25000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 25000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25766

$ws.Range("H85").Value = 25000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 25000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27652

$ws.Range("H102").Value = 1954.8572
$ws.Range("I102").Value = 1781.55
$ws.Range("J102").Value = 2185.9333
$ws.Range("K102").Value = 1781.55
$ws.Range("L102").Value = 2185.9333
$ws.Range("M102").Value = -159.55
$ws.Range("N102").Value = -5429.933300000001

$ws.Range("H123").Value = 10795.6
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 10795.6
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 10795.6
$ws.Range("N123").Value = -15695.6

$ws.Range("H126").Value = 2562.2083
$ws.Range("I126").Value = 1988
$ws.Range("J126").Value = 2906.7334
$ws.Range("K126").Value = 5964
$ws.Range("L126").Value = 8720.200199999999
$ws.Range("M126").Value = -3494
$ws.Range("N126").Value = -13660.2002

$ws.Range("H132").Value = 5456.615
$ws.Range("I132").Value = 4094.8
$ws.Range("J132").Value = 9996
$ws.Range("K132").Value = 12284.4
$ws.Range("L132").Value = 29988
$ws.Range("M132").Value = -9754.400000000001
$ws.Range("N132").Value = -35048

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8992.308000000001
$ws.Range("I22").Value = 1000
$ws.Range("J22").Value = 12544.444
$ws.Range("K22").Value = 1000
$ws.Range("L22").Value = 12544.444
$ws.Range("M22").Value = -705
$ws.Range("N22").Value = -13134.444

$ws.Range("H27").Value = 8992.308000000001
$ws.Range("I27").Value = 1000
$ws.Range("J27").Value = 12544.444
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 12544.444
$ws.Range("M27").Value = -893
$ws.Range("N27").Value = -12758.444

$ws.Range("H55").Value = 470.13333
$ws.Range("I55").Value = 466.66666
$ws.Range("J55").Value = 471
$ws.Range("K55").Value = 466.66666
$ws.Range("L55").Value = 471
$ws.Range("M55").Value = -293.66666
$ws.Range("N55").Value = -817

$ws.Range("H136").Value = 4173.92
$ws.Range("I136").Value = 2281.2632
$ws.Range("J136").Value = 10167.333
$ws.Range("K136").Value = 6843.7896
$ws.Range("L136").Value = 30501.999
$ws.Range("M136").Value = -4293.7896
$ws.Range("N136").Value = -35601.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1740
$ws.Range("I122").Value = 1550
$ws.Range("J122").Value = 1866.6666
$ws.Range("K122").Value = 4650
$ws.Range("L122").Value = 5599.9998
$ws.Range("M122").Value = -2200

$ws.Range("H123").Value = 33146.45
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 33146.45
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 33146.45
$ws.Range("N123").Value = -42946.45

$ws.Range("H126").Value = 72750.92999999999
$ws.Range("I126").Value = 167434.83
$ws.Range("J126").Value = 1738
$ws.Range("K126").Value = 502304.49
$ws.Range("L126").Value = 5214
$ws.Range("M126").Value = -499834.49
